$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 19 (shifts rows 19-27 down to 20-28),
#    for the new "2508" period entry belonging to JULIO MANUEL MORON BATISTA.
$ws.Range("B18:J18").Copy()
$ws.Rows("19:19").Insert()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Fill in the values for the newly inserted row 19.
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73183038"
$ws.Range("D19").Value = "JULIO MANUEL MORON BATISTA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# 3) Re-sequence the late-payment periods for rows 16 and 18 so the four
#    JULIO MANUEL MORON BATISTA periods read 2505, 2506, 2507, 2508 in order.
$ws.Range("E16").Value = "2505"
$ws.Range("E18").Value = "2507"

# 4) Update the summary figures: total overdue value and period count.
$ws.Range("E11").Value = 356434
$ws.Range("F13").Value = 5
